$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 8929
$ws.Range("E2").Value = -1567
$ws.Range("F2").Value = -1567
$ws.Range("G2").Value = -1524
$ws.Range("H2").Value = -2342
$ws.Range("I2").Value = -2338
$ws.Range("J2").Value = -4
$ws.Range("K2").Value = 8526
$ws.Range("L2").Value = 8282
$ws.Range("M2").Value = 244
$ws.Range("N2").Value = 244
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 2669
$ws.Range("Q2").Value = 151
$ws.Range("R2").Value = 3901
$ws.Range("S2").Value = -4511
$ws.Range("T2").Value = 59
$ws.Range("U2").Value = 91
$ws.Range("V2").Value = 3645
$ws.Range("W2").Value = -17.55
$ws.Range("X2").Value = -26.23
$ws.Range("Y2").Value = -143.17
$ws.Range("Z2").Value = -13.92
$ws.Range("AA2").Value = 3392.18
$ws.Range("AB2").Value = -77.23999999999999
$ws.Range("AC2").Value = -218204
$ws.Range("AD2").Value = -0.06
$ws.Range("AE2").Value = 21656
$ws.Range("AF2").Value = 0.6
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 935349

# Row 3
$ws.Range("D3").Value = 6982
$ws.Range("E3").Value = -356
$ws.Range("F3").Value = -356
$ws.Range("G3").Value = -50
$ws.Range("H3").Value = -715
$ws.Range("I3").Value = -715
$ws.Range("K3").Value = 5988
$ws.Range("L3").Value = 5284
$ws.Range("M3").Value = 703
$ws.Range("N3").Value = 703
$ws.Range("P3").Value = 441
$ws.Range("Q3").Value = 1460
$ws.Range("R3").Value = 840
$ws.Range("S3").Value = -1568
$ws.Range("T3").Value = 21
$ws.Range("U3").Value = 1438
$ws.Range("V3").Value = 251
$ws.Range("W3").Value = -5.1
$ws.Range("X3").Value = -10.24
$ws.Range("Y3").Value = -151.07
$ws.Range("Z3").Value = -9.859999999999999
$ws.Range("AA3").Value = 751.4400000000001
$ws.Range("AB3").Value = 324
$ws.Range("AC3").Value = -27665
$ws.Range("AD3").Value = -0.43
$ws.Range("AE3").Value = 16334
$ws.Range("AF3").Value = 0.73
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 4082836
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 5855
$ws.Range("E4").Value = 161
$ws.Range("F4").Value = 161
$ws.Range("G4").Value = 527
$ws.Range("H4").Value = 537
$ws.Range("I4").Value = 537
$ws.Range("K4").Value = 6821
$ws.Range("L4").Value = 4291
$ws.Range("M4").Value = 2531
$ws.Range("N4").Value = 2531
$ws.Range("P4").Value = 1051
$ws.Range("Q4").Value = -100
$ws.Range("R4").Value = 172
$ws.Range("S4").Value = 642
$ws.Range("T4").Value = 4
$ws.Range("U4").Value = -104
$ws.Range("V4").Value = 910
$ws.Range("W4").Value = 2.75
$ws.Range("X4").Value = 9.17
$ws.Range("Y4").Value = 33.2
$ws.Range("Z4").Value = 8.380000000000001
$ws.Range("AA4").Value = 169.54
$ws.Range("AB4").Value = 222.13
$ws.Range("AC4").Value = 6275
$ws.Range("AD4").Value = 1.59
$ws.Range("AE4").Value = 12028
$ws.Range("AF4").Value = 0.83
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 20818499
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 7015
$ws.Range("E5").Value = 256
$ws.Range("F5").Value = 256
$ws.Range("G5").Value = 1119
$ws.Range("H5").Value = 963
$ws.Range("I5").Value = 963
$ws.Range("K5").Value = 6316
$ws.Range("L5").Value = 3275
$ws.Range("M5").Value = 3041
$ws.Range("N5").Value = 3041
$ws.Range("P5").Value = 1054
$ws.Range("Q5").Value = -85
$ws.Range("R5").Value = 604
$ws.Range("S5").Value = -900
$ws.Range("T5").Value = 7
$ws.Range("U5").Value = -92
$ws.Range("V5").Value = 73
$ws.Range("W5").Value = 3.65
$ws.Range("X5").Value = 13.73
$ws.Range("Y5").Value = 34.56
$ws.Range("Z5").Value = 14.66
$ws.Range("AA5").Value = 107.71
$ws.Range("AB5").Value = 313.93
$ws.Range("AC5").Value = 4571
$ws.Range("AD5").Value = 2.21
$ws.Range("AE5").Value = 14429
$ws.Range("AF5").Value = 0.7
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 20853553
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 8982
$ws.Range("E6").Value = 318
$ws.Range("F6").Value = 318
$ws.Range("G6").Value = 590
$ws.Range("H6").Value = 739
$ws.Range("I6").Value = 739
$ws.Range("K6").Value = 7138
$ws.Range("L6").Value = 3519
$ws.Range("M6").Value = 3619
$ws.Range("N6").Value = 3619
$ws.Range("P6").Value = 1055
$ws.Range("Q6").Value = 589
$ws.Range("R6").Value = -1229
$ws.Range("S6").Value = -66
$ws.Range("T6").Value = 5
$ws.Range("U6").Value = 584
$ws.Range("V6").Value = 67
$ws.Range("W6").Value = 3.54
$ws.Range("X6").Value = 8.23
$ws.Range("Y6").Value = 22.21
$ws.Range("Z6").Value = 10.99
$ws.Range("AA6").Value = 97.23
$ws.Range("AB6").Value = 382.6
$ws.Range("AC6").Value = 3507
$ws.Range("AD6").Value = 2.06
$ws.Range("AE6").Value = 17161
$ws.Range("AF6").Value = 0.42
$ws.Range("AG6").Value = 300
$ws.Range("AH6").Value = 4.15
$ws.Range("AI6").Value = 8.57
$ws.Range("AJ6").Value = 20867152

# Row 7
$ws.Range("D7").Value = 10790
$ws.Range("E7").Value = 478
$ws.Range("G7").Value = 598
$ws.Range("H7").Value = 499
$ws.Range("I7").Value = 459
$ws.Range("K7").Value = 8120
$ws.Range("L7").Value = 4029
$ws.Range("M7").Value = 4090
$ws.Range("N7").Value = 4066
$ws.Range("P7").Value = 1082
$ws.Range("Q7").Value = 552
$ws.Range("R7").Value = -502
$ws.Range("S7").Value = -76
$ws.Range("T7").Value = 282
$ws.Range("W7").Value = 4.43
$ws.Range("X7").Value = 4.63
$ws.Range("Y7").Value = 11.95
$ws.Range("Z7").Value = 6.54
$ws.Range("AA7").Value = 98.5
$ws.Range("AC7").Value = 2151
$ws.Range("AD7").Value = 3.86
$ws.Range("AE7").Value = 18442
$ws.Range("AF7").Value = 0.45
$ws.Range("AG7").Value = 300
$ws.Range("AH7").Value = 3.61
$ws.Range("AI7").Value = 14.26
$ws.Range("U7").ClearContents()

# Row 8
$ws.Range("D8").Value = 11677
$ws.Range("E8").Value = 592
$ws.Range("G8").Value = 704
$ws.Range("H8").Value = 617
$ws.Range("I8").Value = 534
$ws.Range("K8").Value = 8227
$ws.Range("L8").Value = 3591
$ws.Range("M8").Value = 4636
$ws.Range("N8").Value = 4611
$ws.Range("P8").Value = 1082
$ws.Range("Q8").Value = 682
$ws.Range("R8").Value = -702
$ws.Range("S8").Value = -133
$ws.Range("T8").Value = 246
$ws.Range("W8").Value = 5.07
$ws.Range("X8").Value = 5.28
$ws.Range("Y8").Value = 12.32
$ws.Range("Z8").Value = 7.55
$ws.Range("AA8").Value = 77.45999999999999
$ws.Range("AC8").Value = 2423
$ws.Range("AD8").Value = 3.46
$ws.Range("AE8").Value = 20912
$ws.Range("AF8").Value = 0.4
$ws.Range("AG8").Value = 350
$ws.Range("AH8").Value = 4.18
$ws.Range("AI8").Value = 14.29
$ws.Range("U8").ClearContents()

# Row 9
$ws.Range("D9").Value = 13806
$ws.Range("E9").Value = 753
$ws.Range("G9").Value = 849
$ws.Range("H9").Value = 641
$ws.Range("I9").Value = 641
$ws.Range("K9").Value = 8734
$ws.Range("L9").Value = 3534
$ws.Range("M9").Value = 5200
$ws.Range("N9").Value = 5175
$ws.Range("P9").Value = 1082
$ws.Range("Q9").Value = 830
$ws.Range("R9").Value = -688
$ws.Range("S9").Value = -103
$ws.Range("T9").Value = 172
$ws.Range("W9").Value = 5.45
$ws.Range("X9").Value = 4.64
$ws.Range("Y9").Value = 13.1
$ws.Range("Z9").Value = 7.56
$ws.Range("AA9").Value = 67.95
$ws.Range("AC9").Value = 2906
$ws.Range("AD9").Value = 2.88
$ws.Range("AE9").Value = 23470
$ws.Range("AF9").Value = 0.36
$ws.Range("AG9").Value = 400
$ws.Range("AH9").Value = 4.77
$ws.Range("AI9").Value = 13.62
$ws.Range("U9").ClearContents()
